# Update the "want-to-go" head-count figures (column F) on worksheet index 1
# (exhibition sheet) and worksheet index 4 (combined/all-types sheet), matching
# the regenerated scrape output.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 613
$ws1.Range("F5").Value = 35
$ws1.Range("F7").Value = 5508
$ws1.Range("F8").Value = 1543
$ws1.Range("F9").Value = 159
$ws1.Range("F10").Value = 3145
$ws1.Range("F13").Value = 1305
$ws1.Range("F14").Value = 4370
$ws1.Range("F15").Value = 1041
$ws1.Range("F17").Value = 1680
$ws1.Range("F18").Value = 2608
$ws1.Range("F20").Value = 38
$ws1.Range("F21").Value = 148
$ws1.Range("F23").Value = 1001
$ws1.Range("F24").Value = 300
$ws1.Range("F29").Value = 1102
$ws1.Range("F30").Value = 394
$ws1.Range("F31").Value = 63
$ws1.Range("F32").Value = 185
$ws1.Range("F33").Value = 320
$ws1.Range("F36").Value = 1702
$ws1.Range("F37").Value = 2213
$ws1.Range("F38").Value = 1030
$ws1.Range("F42").Value = 319
$ws1.Range("F43").Value = 10
$ws1.Range("F46").Value = 411
$ws1.Range("F47").Value = 354
$ws1.Range("F48").Value = 215

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 613
$ws4.Range("F5").Value = 35
$ws4.Range("F7").Value = 5508
$ws4.Range("F8").Value = 1543
$ws4.Range("F9").Value = 159
$ws4.Range("F10").Value = 3145
$ws4.Range("F12").Value = 1305
$ws4.Range("F13").Value = 4370
$ws4.Range("F14").Value = 1041
$ws4.Range("F15").Value = 1680
$ws4.Range("F20").Value = 38
$ws4.Range("F21").Value = 148
$ws4.Range("F24").Value = 1001
$ws4.Range("F25").Value = 300
$ws4.Range("F30").Value = 1102
$ws4.Range("F31").Value = 394
$ws4.Range("F32").Value = 63
$ws4.Range("F33").Value = 185
$ws4.Range("F36").Value = 1702
$ws4.Range("F37").Value = 2213
$ws4.Range("F38").Value = 1030
$ws4.Range("F43").Value = 319
$ws4.Range("F45").Value = 411
$ws4.Range("F46").Value = 354
$ws4.Range("F47").Value = 215

